# Auto-generated Excel COM-interop script to apply scheduled market-price/profit updates
# to the Ravana_Profits leve-crafting workbook (chore: update Sheets via scheduled runner).
#
# For every changed Leve row we overwrite the recomputed market-price / profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ).
# Some rows additionally gain or lose a LeveProfitNQ (M) / LeveProfitHQ (N) cell entirely
# (e.g. when NQ or HQ crafting stops being profitable/possible) - those are cleared with
# ClearContents() so the cell is removed from the sheet rather than merely zeroed.

$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 29999.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 29999.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 29999.5
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -30967.5
# Row 107
$ws.Range("H107").Value = 238.66667
$ws.Range("I107").Value = 256.125
$ws.Range("J107").Value = 99
$ws.Range("K107").Value = 256.125
$ws.Range("L107").Value = 99
$ws.Range("M107").Value = 1663.875
$ws.Range("N107").Value = -3939
# Row 116
$ws.Range("H116").Value = 5099.1113
$ws.Range("I116").Value = 4999
$ws.Range("J116").Value = 5299.3335
$ws.Range("K116").Value = 4999
$ws.Range("L116").Value = 5299.3335
$ws.Range("M116").Value = -1557
$ws.Range("N116").Value = -12183.3335
# Row 132
$ws.Range("H132").Value = 1376.5385
$ws.Range("I132").Value = 1335.3636
$ws.Range("K132").Value = 4006.0908
$ws.Range("M132").Value = -1476.0908
# Row 137
$ws.Range("H137").Value = 2617.2144
$ws.Range("I137").Value = 1908.8572
$ws.Range("J137").Value = 4742.2856
$ws.Range("K137").Value = 5726.571599999999
$ws.Range("L137").Value = 14226.8568
$ws.Range("M137").Value = -3176.571599999999
$ws.Range("N137").Value = -19326.8568
# Row 138
$ws.Range("H138").Value = 3450.8
$ws.Range("I138").Value = 1429.5264
$ws.Range("J138").Value = 5851.0625
$ws.Range("K138").Value = 4288.5792
$ws.Range("L138").Value = 17553.1875
$ws.Range("M138").Value = 851.4207999999999
$ws.Range("N138").Value = -27833.1875

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3702.9707
$ws.Range("I32").Value = 2964.9033
$ws.Range("K32").Value = 2964.9033
$ws.Range("M32").Value = -2677.9033
# Row 45
$ws.Range("H45").Value = 1983.4117
$ws.Range("I45").Value = 1826.1875
$ws.Range("K45").Value = 1826.1875
$ws.Range("M45").Value = -1449.1875
# Row 61
$ws.Range("H61").Value = 2400.6365
$ws.Range("I61").Value = 2461
$ws.Range("J61").Value = 2239.6667
$ws.Range("K61").Value = 2461
$ws.Range("L61").Value = 2239.6667
$ws.Range("M61").Value = -2249
$ws.Range("N61").Value = -2663.6667
# Row 74
$ws.Range("H74").Value = 2418.08
$ws.Range("I74").Value = 2526.2173
$ws.Range("K74").Value = 2526.2173
$ws.Range("M74").Value = -1652.2173
# Row 77
$ws.Range("H77").Value = 2418.08
$ws.Range("I77").Value = 2526.2173
$ws.Range("K77").Value = 12631.0865
$ws.Range("M77").Value = -8263.086499999999
# Row 97
$ws.Range("H97").Value = 799.0625
$ws.Range("I97").Value = 381.5
$ws.Range("J97").Value = 3722
$ws.Range("K97").Value = 381.5
$ws.Range("L97").Value = 3722
$ws.Range("M97").Value = 114.5
$ws.Range("N97").Value = -4714
# Row 132
$ws.Range("H132").Value = 2076.8462
$ws.Range("I132").Value = 1735.9688
$ws.Range("J132").Value = 3635.1428
$ws.Range("K132").Value = 5207.9064
$ws.Range("L132").Value = 10905.4284
$ws.Range("M132").Value = -2677.9064
$ws.Range("N132").Value = -15965.4284
# Row 136
$ws.Range("H136").Value = 2400.6365
$ws.Range("I136").Value = 2461
$ws.Range("J136").Value = 2239.6667
$ws.Range("K136").Value = 7383
$ws.Range("L136").Value = 6719.000100000001
$ws.Range("M136").Value = -4833
$ws.Range("N136").Value = -11819.0001

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2352.742
$ws.Range("I134").Value = 2345.08
$ws.Range("K134").Value = 7035.24
$ws.Range("M134").Value = -4500.24

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1459.2858
$ws.Range("I16").Value = 1243.2
$ws.Range("K16").Value = 1243.2
$ws.Range("M16").Value = -956.2
# Row 22
$ws.Range("H22").Value = 984.3333
$ws.Range("I22").Value = 976
$ws.Range("K22").Value = 976
$ws.Range("M22").Value = -626
# Row 113
$ws.Range("H113").Value = 1459.2858
$ws.Range("I113").Value = 1243.2
$ws.Range("K113").Value = 1243.2
$ws.Range("M113").Value = 926.8
# Row 131
$ws.Range("H131").Value = 95000
$ws.Range("J131").Value = 95000
$ws.Range("L131").Value = 95000
$ws.Range("N131").Value = -105080
# Row 134
$ws.Range("H134").Value = 4336.1875
$ws.Range("I134").Value = 4822.1665
$ws.Range("K134").Value = 14466.4995
$ws.Range("M134").Value = -11931.4995

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 130.5
$ws.Range("J12").Value = 90.111115
$ws.Range("L12").Value = 270.333345
$ws.Range("N12").Value = -616.333345
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 122
$ws.Range("H122").Value = 1100
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9900
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -14800
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 122
$ws.Range("H122").Value = 688
$ws.Range("I122").Value = 674.1667
$ws.Range("J122").Value = 771
$ws.Range("K122").Value = 2022.5001
$ws.Range("L122").Value = 2313
$ws.Range("M122").Value = 427.4999
$ws.Range("N122").Value = -7213
# Row 132
$ws.Range("H132").Value = 3074.4285
$ws.Range("I132").Value = 2731.3635
$ws.Range("K132").Value = 8194.0905
$ws.Range("M132").Value = -5664.0905

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 273.5
$ws.Range("I55").Value = 248
$ws.Range("K55").Value = 248
$ws.Range("M55").Value = -75
# Row 132
$ws.Range("H132").Value = 3830.8333
$ws.Range("J132").Value = 4497.5
$ws.Range("L132").Value = 13492.5
$ws.Range("N132").Value = -18552.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2568.0908
$ws.Range("I132").Value = 1786.0714
$ws.Range("J132").Value = 3936.625
$ws.Range("K132").Value = 5358.2142
$ws.Range("L132").Value = 11809.875
$ws.Range("M132").Value = -2828.2142
$ws.Range("N132").Value = -16869.875
# Row 137
$ws.Range("H137").Value = 53000
$ws.Range("J137").Value = 53000
$ws.Range("L137").Value = 53000
$ws.Range("N137").Value = -63200

